$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window view: windowHeight 10990 -> 12490 (reflected in the saved workbook view)
$excel.Windows.Item(1).Height = 12490

# Remove trailing "*" from these country name pairs (EN/FR columns B/C)
$ws.Range("B83").Value = "Botswana"
$ws.Range("C83").Value = "Botswana"

$ws.Range("B87").Value = "Chile"
$ws.Range("C87").Value = "Chili"

$ws.Range("B91").Value = "DR Congo"
$ws.Range("C91").Value = "RD Congo"

$ws.Range("B119").Value = "Ghana"
$ws.Range("C119").Value = "Ghana"

$ws.Range("B127").Value = "Guyana"
$ws.Range("C127").Value = "Guyana"

$ws.Range("B195").Value = "Papua New Guinea"
$ws.Range("C195").Value = "Papouasie-Nouvelle-Guinée"

$ws.Range("B214").Value = "South Sudan"
$ws.Range("C214").Value = "Soudan du Sud"

$ws.Range("B224").Value = "Togo"
$ws.Range("C224").Value = "Togo"

$ws.Range("B247").Value = "Zambia"
$ws.Range("C247").Value = "Zambie"

# Add trailing "*" to Timor-Leste entries
$ws.Range("B228").Value = "Timor-Leste*"
$ws.Range("C228").Value = "Timor-Leste*"
